$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (s=1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J (mirror each other row by row)
$values = @{
    2  = 6
    3  = 8
    4  = 8
    5  = 8
    6  = 8
    7  = 7
    8  = 7
    9  = 5
    10 = 8
    11 = 4
    12 = 4
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 9).Value = $val
    $ws.Cells.Item($row, 10).Value = $val
}
